# Generate Report for Handoff
#
# This regenerates the latest-handoff timestamps for the file
# "90d64a18-4662-463c-af3a-091e1c5bc0bd" row (row 5) across the
# Overview sheet and each per-locale handoff sheet, as part of
# producing a fresh localization-status report.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest Handoff Date" column (D) for the
# 90d64a18-4662-463c-af3a-091e1c5bc0bd.md row.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("D5").Value = "2016-44-13 18:44:55"

# zh-cn sheet: "Latest Handoff Datetime" column (E) for the
# 90d64a18-4662-463c-af3a-091e1c5bc0bd row.
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E5").Value = "2016-03-13 18:44:46"

# de-de sheet: "Latest Handoff Datetime" column (E) for the
# 90d64a18-4662-463c-af3a-091e1c5bc0bd row.
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E5").Value = "2016-03-13 18:44:55"
